# Rename the header row from the generic "_old"/"_new" suffixes to the
# concrete format-version suffixes "_FV2210"/"_FV2304", then turn the
# header + data range into a real Excel Table (ListObject) and freeze
# the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells A1:U1 -------------------------------------
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn A1:U76 into an Excel Table (adds xl/tables/table1.xml,
#        autofilter + tableParts reference) -----------------------------
$range = $ws.Range("A1:U76")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
